$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Names: ', ' joins between co-contractors -> '. ' (8 strings / 10 cells) ---
$ws.Range("E31").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E181").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E87").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E106").Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("F106").Value = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
$ws.Range("E116").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E160").Value = "DODERA. JORGE ABELARDO"
$ws.Range("E166").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("F126").Value = "MERCANZINI. GASTON ARIEL"
$ws.Range("F137").Value = "OLVEIRA. ALBERTO MIGUEL"

# --- Importe: es-AR '1.234,56' -> plain '1234.56' text (213 cells) ---
# NumberFormat '@' forces the assignment to stay text instead of being
# parsed into a real number; ClearFormats() afterwards drops the cell-level
# style override again so only the *value* changes, matching the source diff.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "27920.00"
$ws.Range("H2").ClearFormats()
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "3800.00"
$ws.Range("H3").ClearFormats()
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "50000.00"
$ws.Range("H4").ClearFormats()
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "0.50"
$ws.Range("H5").ClearFormats()
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "3665.00"
$ws.Range("H6").ClearFormats()
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "257282.40"
$ws.Range("H7").ClearFormats()
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "1620.00"
$ws.Range("H8").ClearFormats()
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1103.40"
$ws.Range("H9").ClearFormats()
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "310.00"
$ws.Range("H10").ClearFormats()
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "87801.76"
$ws.Range("H11").ClearFormats()
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "2334.00"
$ws.Range("H12").ClearFormats()
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "52307.95"
$ws.Range("H13").ClearFormats()
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "2311.20"
$ws.Range("H14").ClearFormats()
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "6475.71"
$ws.Range("H15").ClearFormats()
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "4905.00"
$ws.Range("H16").ClearFormats()
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "1450.00"
$ws.Range("H17").ClearFormats()
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "250.00"
$ws.Range("H18").ClearFormats()
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "1200.00"
$ws.Range("H19").ClearFormats()
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "1760.00"
$ws.Range("H20").ClearFormats()
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "100.48"
$ws.Range("H21").ClearFormats()
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "39.40"
$ws.Range("H22").ClearFormats()
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "16.50"
$ws.Range("H23").ClearFormats()
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "640.00"
$ws.Range("H24").ClearFormats()
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "163865.73"
$ws.Range("H25").ClearFormats()
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "8400.00"
$ws.Range("H26").ClearFormats()
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "583.86"
$ws.Range("H27").ClearFormats()
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "692.00"
$ws.Range("H28").ClearFormats()
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "230707.80"
$ws.Range("H29").ClearFormats()
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "13.26"
$ws.Range("H30").ClearFormats()
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "380.00"
$ws.Range("H31").ClearFormats()
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "100.00"
$ws.Range("H32").ClearFormats()
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "60.00"
$ws.Range("H33").ClearFormats()
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "25.56"
$ws.Range("H34").ClearFormats()
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "2102.69"
$ws.Range("H35").ClearFormats()
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "61391.00"
$ws.Range("H36").ClearFormats()
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "1547.00"
$ws.Range("H37").ClearFormats()
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "128681.35"
$ws.Range("H38").ClearFormats()
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "8.00"
$ws.Range("H39").ClearFormats()
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "620.00"
$ws.Range("H40").ClearFormats()
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "48.00"
$ws.Range("H41").ClearFormats()
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "3094.26"
$ws.Range("H42").ClearFormats()
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "40.00"
$ws.Range("H43").ClearFormats()
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "1960.00"
$ws.Range("H44").ClearFormats()
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "751.18"
$ws.Range("H45").ClearFormats()
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "3141.60"
$ws.Range("H46").ClearFormats()
$ws.Range("H47").NumberFormat = "@"
$ws.Range("H47").Value = "15957.00"
$ws.Range("H47").ClearFormats()
$ws.Range("H48").NumberFormat = "@"
$ws.Range("H48").Value = "24.00"
$ws.Range("H48").ClearFormats()
$ws.Range("H49").NumberFormat = "@"
$ws.Range("H49").Value = "1596.90"
$ws.Range("H49").ClearFormats()
$ws.Range("H50").NumberFormat = "@"
$ws.Range("H50").Value = "10737.54"
$ws.Range("H50").ClearFormats()
$ws.Range("H51").NumberFormat = "@"
$ws.Range("H51").Value = "314.47"
$ws.Range("H51").ClearFormats()
$ws.Range("H52").NumberFormat = "@"
$ws.Range("H52").Value = "170.00"
$ws.Range("H52").ClearFormats()
$ws.Range("H53").NumberFormat = "@"
$ws.Range("H53").Value = "4849.08"
$ws.Range("H53").ClearFormats()
$ws.Range("H54").NumberFormat = "@"
$ws.Range("H54").Value = "28750.77"
$ws.Range("H54").ClearFormats()
$ws.Range("H55").NumberFormat = "@"
$ws.Range("H55").Value = "7725.75"
$ws.Range("H55").ClearFormats()
$ws.Range("H56").NumberFormat = "@"
$ws.Range("H56").Value = "1835.00"
$ws.Range("H56").ClearFormats()
$ws.Range("H57").NumberFormat = "@"
$ws.Range("H57").Value = "31984.60"
$ws.Range("H57").ClearFormats()
$ws.Range("H58").NumberFormat = "@"
$ws.Range("H58").Value = "2670.36"
$ws.Range("H58").ClearFormats()
$ws.Range("H59").NumberFormat = "@"
$ws.Range("H59").Value = "3809.08"
$ws.Range("H59").ClearFormats()
$ws.Range("H60").NumberFormat = "@"
$ws.Range("H60").Value = "2974.00"
$ws.Range("H60").ClearFormats()
$ws.Range("H61").NumberFormat = "@"
$ws.Range("H61").Value = "557.15"
$ws.Range("H61").ClearFormats()
$ws.Range("H62").NumberFormat = "@"
$ws.Range("H62").Value = "534.99"
$ws.Range("H62").ClearFormats()
$ws.Range("H63").NumberFormat = "@"
$ws.Range("H63").Value = "21670.00"
$ws.Range("H63").ClearFormats()
$ws.Range("H64").NumberFormat = "@"
$ws.Range("H64").Value = "3736.56"
$ws.Range("H64").ClearFormats()
$ws.Range("H65").NumberFormat = "@"
$ws.Range("H65").Value = "2753.00"
$ws.Range("H65").ClearFormats()
$ws.Range("H66").NumberFormat = "@"
$ws.Range("H66").Value = "4968.82"
$ws.Range("H66").ClearFormats()
$ws.Range("H67").NumberFormat = "@"
$ws.Range("H67").Value = "28000.00"
$ws.Range("H67").ClearFormats()
$ws.Range("H68").NumberFormat = "@"
$ws.Range("H68").Value = "1628.66"
$ws.Range("H68").ClearFormats()
$ws.Range("H69").NumberFormat = "@"
$ws.Range("H69").Value = "118.65"
$ws.Range("H69").ClearFormats()
$ws.Range("H70").NumberFormat = "@"
$ws.Range("H70").Value = "2240.00"
$ws.Range("H70").ClearFormats()
$ws.Range("H71").NumberFormat = "@"
$ws.Range("H71").Value = "42000.00"
$ws.Range("H71").ClearFormats()
$ws.Range("H72").NumberFormat = "@"
$ws.Range("H72").Value = "6691.00"
$ws.Range("H72").ClearFormats()
$ws.Range("H73").NumberFormat = "@"
$ws.Range("H73").Value = "4445.00"
$ws.Range("H73").ClearFormats()
$ws.Range("H74").NumberFormat = "@"
$ws.Range("H74").Value = "56840.00"
$ws.Range("H74").ClearFormats()
$ws.Range("H75").NumberFormat = "@"
$ws.Range("H75").Value = "3850.00"
$ws.Range("H75").ClearFormats()
$ws.Range("H76").NumberFormat = "@"
$ws.Range("H76").Value = "25560.00"
$ws.Range("H76").ClearFormats()
$ws.Range("H77").NumberFormat = "@"
$ws.Range("H77").Value = "380.00"
$ws.Range("H77").ClearFormats()
$ws.Range("H78").NumberFormat = "@"
$ws.Range("H78").Value = "14040.00"
$ws.Range("H78").ClearFormats()
$ws.Range("H79").NumberFormat = "@"
$ws.Range("H79").Value = "2660.00"
$ws.Range("H79").ClearFormats()
$ws.Range("H80").NumberFormat = "@"
$ws.Range("H80").Value = "22796.00"
$ws.Range("H80").ClearFormats()
$ws.Range("H81").NumberFormat = "@"
$ws.Range("H81").Value = "1380.00"
$ws.Range("H81").ClearFormats()
$ws.Range("H82").NumberFormat = "@"
$ws.Range("H82").Value = "21196.80"
$ws.Range("H82").ClearFormats()
$ws.Range("H83").NumberFormat = "@"
$ws.Range("H83").Value = "6000.00"
$ws.Range("H83").ClearFormats()
$ws.Range("H84").NumberFormat = "@"
$ws.Range("H84").Value = "3000.00"
$ws.Range("H84").ClearFormats()
$ws.Range("H85").NumberFormat = "@"
$ws.Range("H85").Value = "289.00"
$ws.Range("H85").ClearFormats()
$ws.Range("H86").NumberFormat = "@"
$ws.Range("H86").Value = "7587.00"
$ws.Range("H86").ClearFormats()
$ws.Range("H87").NumberFormat = "@"
$ws.Range("H87").Value = "1820.00"
$ws.Range("H87").ClearFormats()
$ws.Range("H88").NumberFormat = "@"
$ws.Range("H88").Value = "20916.00"
$ws.Range("H88").ClearFormats()
$ws.Range("H89").NumberFormat = "@"
$ws.Range("H89").Value = "512.70"
$ws.Range("H89").ClearFormats()
$ws.Range("H90").NumberFormat = "@"
$ws.Range("H90").Value = "16980.00"
$ws.Range("H90").ClearFormats()
$ws.Range("H91").NumberFormat = "@"
$ws.Range("H91").Value = "2600.00"
$ws.Range("H91").ClearFormats()
$ws.Range("H92").NumberFormat = "@"
$ws.Range("H92").Value = "270000.00"
$ws.Range("H92").ClearFormats()
$ws.Range("H93").NumberFormat = "@"
$ws.Range("H93").Value = "154320.00"
$ws.Range("H93").ClearFormats()
$ws.Range("H94").NumberFormat = "@"
$ws.Range("H94").Value = "334375.87"
$ws.Range("H94").ClearFormats()
$ws.Range("H95").NumberFormat = "@"
$ws.Range("H95").Value = "139500.00"
$ws.Range("H95").ClearFormats()
$ws.Range("H96").NumberFormat = "@"
$ws.Range("H96").Value = "5846.78"
$ws.Range("H96").ClearFormats()
$ws.Range("H97").NumberFormat = "@"
$ws.Range("H97").Value = "1320.00"
$ws.Range("H97").ClearFormats()
$ws.Range("H98").NumberFormat = "@"
$ws.Range("H98").Value = "9.25"
$ws.Range("H98").ClearFormats()
$ws.Range("H99").NumberFormat = "@"
$ws.Range("H99").Value = "17.78"
$ws.Range("H99").ClearFormats()
$ws.Range("H100").NumberFormat = "@"
$ws.Range("H100").Value = "699.66"
$ws.Range("H100").ClearFormats()
$ws.Range("H101").NumberFormat = "@"
$ws.Range("H101").Value = "8905.00"
$ws.Range("H101").ClearFormats()
$ws.Range("H102").NumberFormat = "@"
$ws.Range("H102").Value = "490.50"
$ws.Range("H102").ClearFormats()
$ws.Range("H103").NumberFormat = "@"
$ws.Range("H103").Value = "360.00"
$ws.Range("H103").ClearFormats()
$ws.Range("H104").NumberFormat = "@"
$ws.Range("H104").Value = "2512.00"
$ws.Range("H104").ClearFormats()
$ws.Range("H105").NumberFormat = "@"
$ws.Range("H105").Value = "240.00"
$ws.Range("H105").ClearFormats()
$ws.Range("H106").NumberFormat = "@"
$ws.Range("H106").Value = "250.00"
$ws.Range("H106").ClearFormats()
$ws.Range("H107").NumberFormat = "@"
$ws.Range("H107").Value = "5793.65"
$ws.Range("H107").ClearFormats()
$ws.Range("H108").NumberFormat = "@"
$ws.Range("H108").Value = "311.20"
$ws.Range("H108").ClearFormats()
$ws.Range("H109").NumberFormat = "@"
$ws.Range("H109").Value = "297.00"
$ws.Range("H109").ClearFormats()
$ws.Range("H110").NumberFormat = "@"
$ws.Range("H110").Value = "228.68"
$ws.Range("H110").ClearFormats()
$ws.Range("H111").NumberFormat = "@"
$ws.Range("H111").Value = "97.10"
$ws.Range("H111").ClearFormats()
$ws.Range("H112").NumberFormat = "@"
$ws.Range("H112").Value = "1518.00"
$ws.Range("H112").ClearFormats()
$ws.Range("H113").NumberFormat = "@"
$ws.Range("H113").Value = "3747.08"
$ws.Range("H113").ClearFormats()
$ws.Range("H114").NumberFormat = "@"
$ws.Range("H114").Value = "6702.71"
$ws.Range("H114").ClearFormats()
$ws.Range("H115").NumberFormat = "@"
$ws.Range("H115").Value = "3300.00"
$ws.Range("H115").ClearFormats()
$ws.Range("H116").NumberFormat = "@"
$ws.Range("H116").Value = "1000.00"
$ws.Range("H116").ClearFormats()
$ws.Range("H117").NumberFormat = "@"
$ws.Range("H117").Value = "1842.00"
$ws.Range("H117").ClearFormats()
$ws.Range("H118").NumberFormat = "@"
$ws.Range("H118").Value = "30080.00"
$ws.Range("H118").ClearFormats()
$ws.Range("H119").NumberFormat = "@"
$ws.Range("H119").Value = "3980.00"
$ws.Range("H119").ClearFormats()
$ws.Range("H120").NumberFormat = "@"
$ws.Range("H120").Value = "1980.00"
$ws.Range("H120").ClearFormats()
$ws.Range("H121").NumberFormat = "@"
$ws.Range("H121").Value = "23100.00"
$ws.Range("H121").ClearFormats()
$ws.Range("H122").NumberFormat = "@"
$ws.Range("H122").Value = "132000.00"
$ws.Range("H122").ClearFormats()
$ws.Range("H123").NumberFormat = "@"
$ws.Range("H123").Value = "3150.00"
$ws.Range("H123").ClearFormats()
$ws.Range("H124").NumberFormat = "@"
$ws.Range("H124").Value = "800.00"
$ws.Range("H124").ClearFormats()
$ws.Range("H125").NumberFormat = "@"
$ws.Range("H125").Value = "3800.00"
$ws.Range("H125").ClearFormats()
$ws.Range("H126").NumberFormat = "@"
$ws.Range("H126").Value = "6000.00"
$ws.Range("H126").ClearFormats()
$ws.Range("H127").NumberFormat = "@"
$ws.Range("H127").Value = "298685.00"
$ws.Range("H127").ClearFormats()
$ws.Range("H128").NumberFormat = "@"
$ws.Range("H128").Value = "1276.00"
$ws.Range("H128").ClearFormats()
$ws.Range("H129").NumberFormat = "@"
$ws.Range("H129").Value = "66.51"
$ws.Range("H129").ClearFormats()
$ws.Range("H130").NumberFormat = "@"
$ws.Range("H130").Value = "1556.00"
$ws.Range("H130").ClearFormats()
$ws.Range("H131").NumberFormat = "@"
$ws.Range("H131").Value = "51.76"
$ws.Range("H131").ClearFormats()
$ws.Range("H132").NumberFormat = "@"
$ws.Range("H132").Value = "169.66"
$ws.Range("H132").ClearFormats()
$ws.Range("H133").NumberFormat = "@"
$ws.Range("H133").Value = "20.00"
$ws.Range("H133").ClearFormats()
$ws.Range("H134").NumberFormat = "@"
$ws.Range("H134").Value = "1734.00"
$ws.Range("H134").ClearFormats()
$ws.Range("H135").NumberFormat = "@"
$ws.Range("H135").Value = "132623.60"
$ws.Range("H135").ClearFormats()
$ws.Range("H136").NumberFormat = "@"
$ws.Range("H136").Value = "83865.60"
$ws.Range("H136").ClearFormats()
$ws.Range("H137").NumberFormat = "@"
$ws.Range("H137").Value = "3025.00"
$ws.Range("H137").ClearFormats()
$ws.Range("H138").NumberFormat = "@"
$ws.Range("H138").Value = "3960.00"
$ws.Range("H138").ClearFormats()
$ws.Range("H139").NumberFormat = "@"
$ws.Range("H139").Value = "500.00"
$ws.Range("H139").ClearFormats()
$ws.Range("H140").NumberFormat = "@"
$ws.Range("H140").Value = "18000.00"
$ws.Range("H140").ClearFormats()
$ws.Range("H141").NumberFormat = "@"
$ws.Range("H141").Value = "500.00"
$ws.Range("H141").ClearFormats()
$ws.Range("H142").NumberFormat = "@"
$ws.Range("H142").Value = "7114.80"
$ws.Range("H142").ClearFormats()
$ws.Range("H143").NumberFormat = "@"
$ws.Range("H143").Value = "460.00"
$ws.Range("H143").ClearFormats()
$ws.Range("H144").NumberFormat = "@"
$ws.Range("H144").Value = "1600.00"
$ws.Range("H144").ClearFormats()
$ws.Range("H145").NumberFormat = "@"
$ws.Range("H145").Value = "1000.00"
$ws.Range("H145").ClearFormats()
$ws.Range("H146").NumberFormat = "@"
$ws.Range("H146").Value = "3000.00"
$ws.Range("H146").ClearFormats()
$ws.Range("H147").NumberFormat = "@"
$ws.Range("H147").Value = "3315.00"
$ws.Range("H147").ClearFormats()
$ws.Range("H148").NumberFormat = "@"
$ws.Range("H148").Value = "1450.00"
$ws.Range("H148").ClearFormats()
$ws.Range("H149").NumberFormat = "@"
$ws.Range("H149").Value = "700.00"
$ws.Range("H149").ClearFormats()
$ws.Range("H150").NumberFormat = "@"
$ws.Range("H150").Value = "7500.00"
$ws.Range("H150").ClearFormats()
$ws.Range("H151").NumberFormat = "@"
$ws.Range("H151").Value = "4962.00"
$ws.Range("H151").ClearFormats()
$ws.Range("H152").NumberFormat = "@"
$ws.Range("H152").Value = "2800.00"
$ws.Range("H152").ClearFormats()
$ws.Range("H153").NumberFormat = "@"
$ws.Range("H153").Value = "200.00"
$ws.Range("H153").ClearFormats()
$ws.Range("H154").NumberFormat = "@"
$ws.Range("H154").Value = "700.00"
$ws.Range("H154").ClearFormats()
$ws.Range("H155").NumberFormat = "@"
$ws.Range("H155").Value = "678.00"
$ws.Range("H155").ClearFormats()
$ws.Range("H156").NumberFormat = "@"
$ws.Range("H156").Value = "150.00"
$ws.Range("H156").ClearFormats()
$ws.Range("H157").NumberFormat = "@"
$ws.Range("H157").Value = "350.00"
$ws.Range("H157").ClearFormats()
$ws.Range("H158").NumberFormat = "@"
$ws.Range("H158").Value = "290.00"
$ws.Range("H158").ClearFormats()
$ws.Range("H159").NumberFormat = "@"
$ws.Range("H159").Value = "1090.00"
$ws.Range("H159").ClearFormats()
$ws.Range("H160").NumberFormat = "@"
$ws.Range("H160").Value = "800.00"
$ws.Range("H160").ClearFormats()
$ws.Range("H161").NumberFormat = "@"
$ws.Range("H161").Value = "2350.00"
$ws.Range("H161").ClearFormats()
$ws.Range("H162").NumberFormat = "@"
$ws.Range("H162").Value = "480.00"
$ws.Range("H162").ClearFormats()
$ws.Range("H163").NumberFormat = "@"
$ws.Range("H163").Value = "6550.00"
$ws.Range("H163").ClearFormats()
$ws.Range("H164").NumberFormat = "@"
$ws.Range("H164").Value = "443.13"
$ws.Range("H164").ClearFormats()
$ws.Range("H165").NumberFormat = "@"
$ws.Range("H165").Value = "1179.00"
$ws.Range("H165").ClearFormats()
$ws.Range("H166").NumberFormat = "@"
$ws.Range("H166").Value = "280.00"
$ws.Range("H166").ClearFormats()
$ws.Range("H167").NumberFormat = "@"
$ws.Range("H167").Value = "1127.00"
$ws.Range("H167").ClearFormats()
$ws.Range("H168").NumberFormat = "@"
$ws.Range("H168").Value = "1184.74"
$ws.Range("H168").ClearFormats()
$ws.Range("H169").NumberFormat = "@"
$ws.Range("H169").Value = "19433.99"
$ws.Range("H169").ClearFormats()
$ws.Range("H170").NumberFormat = "@"
$ws.Range("H170").Value = "2200.00"
$ws.Range("H170").ClearFormats()
$ws.Range("H171").NumberFormat = "@"
$ws.Range("H171").Value = "452.00"
$ws.Range("H171").ClearFormats()
$ws.Range("H172").NumberFormat = "@"
$ws.Range("H172").Value = "5400.00"
$ws.Range("H172").ClearFormats()
$ws.Range("H173").NumberFormat = "@"
$ws.Range("H173").Value = "4006.55"
$ws.Range("H173").ClearFormats()
$ws.Range("H174").NumberFormat = "@"
$ws.Range("H174").Value = "209.80"
$ws.Range("H174").ClearFormats()
$ws.Range("H175").NumberFormat = "@"
$ws.Range("H175").Value = "404.00"
$ws.Range("H175").ClearFormats()
$ws.Range("H176").NumberFormat = "@"
$ws.Range("H176").Value = "322.00"
$ws.Range("H176").ClearFormats()
$ws.Range("H177").NumberFormat = "@"
$ws.Range("H177").Value = "420.14"
$ws.Range("H177").ClearFormats()
$ws.Range("H178").NumberFormat = "@"
$ws.Range("H178").Value = "220.00"
$ws.Range("H178").ClearFormats()
$ws.Range("H179").NumberFormat = "@"
$ws.Range("H179").Value = "17166.01"
$ws.Range("H179").ClearFormats()
$ws.Range("H180").NumberFormat = "@"
$ws.Range("H180").Value = "18.00"
$ws.Range("H180").ClearFormats()
$ws.Range("H181").NumberFormat = "@"
$ws.Range("H181").Value = "2536.00"
$ws.Range("H181").ClearFormats()
$ws.Range("H182").NumberFormat = "@"
$ws.Range("H182").Value = "657.00"
$ws.Range("H182").ClearFormats()
$ws.Range("H183").NumberFormat = "@"
$ws.Range("H183").Value = "7055.00"
$ws.Range("H183").ClearFormats()
$ws.Range("H184").NumberFormat = "@"
$ws.Range("H184").Value = "84.00"
$ws.Range("H184").ClearFormats()
$ws.Range("H185").NumberFormat = "@"
$ws.Range("H185").Value = "112.77"
$ws.Range("H185").ClearFormats()
$ws.Range("H186").NumberFormat = "@"
$ws.Range("H186").Value = "854.04"
$ws.Range("H186").ClearFormats()
$ws.Range("H187").NumberFormat = "@"
$ws.Range("H187").Value = "24.07"
$ws.Range("H187").ClearFormats()
$ws.Range("H188").NumberFormat = "@"
$ws.Range("H188").Value = "3890.00"
$ws.Range("H188").ClearFormats()
$ws.Range("H189").NumberFormat = "@"
$ws.Range("H189").Value = "185.55"
$ws.Range("H189").ClearFormats()
$ws.Range("H190").NumberFormat = "@"
$ws.Range("H190").Value = "9670.90"
$ws.Range("H190").ClearFormats()
$ws.Range("H191").NumberFormat = "@"
$ws.Range("H191").Value = "40.00"
$ws.Range("H191").ClearFormats()
$ws.Range("H192").NumberFormat = "@"
$ws.Range("H192").Value = "2613.90"
$ws.Range("H192").ClearFormats()
$ws.Range("H193").NumberFormat = "@"
$ws.Range("H193").Value = "6554.33"
$ws.Range("H193").ClearFormats()
$ws.Range("H194").NumberFormat = "@"
$ws.Range("H194").Value = "35910.00"
$ws.Range("H194").ClearFormats()
$ws.Range("H195").NumberFormat = "@"
$ws.Range("H195").Value = "2050.00"
$ws.Range("H195").ClearFormats()
$ws.Range("H196").NumberFormat = "@"
$ws.Range("H196").Value = "28500.00"
$ws.Range("H196").ClearFormats()
$ws.Range("H197").NumberFormat = "@"
$ws.Range("H197").Value = "7900.00"
$ws.Range("H197").ClearFormats()
$ws.Range("H198").NumberFormat = "@"
$ws.Range("H198").Value = "3900.00"
$ws.Range("H198").ClearFormats()
$ws.Range("H199").NumberFormat = "@"
$ws.Range("H199").Value = "112000.00"
$ws.Range("H199").ClearFormats()
$ws.Range("H200").NumberFormat = "@"
$ws.Range("H200").Value = "145000.00"
$ws.Range("H200").ClearFormats()
$ws.Range("H201").NumberFormat = "@"
$ws.Range("H201").Value = "432000.00"
$ws.Range("H201").ClearFormats()
$ws.Range("H202").NumberFormat = "@"
$ws.Range("H202").Value = "1347594.50"
$ws.Range("H202").ClearFormats()
$ws.Range("H203").NumberFormat = "@"
$ws.Range("H203").Value = "390538.69"
$ws.Range("H203").ClearFormats()
$ws.Range("H204").NumberFormat = "@"
$ws.Range("H204").Value = "208565.00"
$ws.Range("H204").ClearFormats()
$ws.Range("H205").NumberFormat = "@"
$ws.Range("H205").Value = "7200.10"
$ws.Range("H205").ClearFormats()
$ws.Range("H206").NumberFormat = "@"
$ws.Range("H206").Value = "81520.00"
$ws.Range("H206").ClearFormats()
$ws.Range("H207").NumberFormat = "@"
$ws.Range("H207").Value = "4143.24"
$ws.Range("H207").ClearFormats()
$ws.Range("H208").NumberFormat = "@"
$ws.Range("H208").Value = "17400.00"
$ws.Range("H208").ClearFormats()
$ws.Range("H209").NumberFormat = "@"
$ws.Range("H209").Value = "4775.00"
$ws.Range("H209").ClearFormats()
$ws.Range("H210").NumberFormat = "@"
$ws.Range("H210").Value = "6400.00"
$ws.Range("H210").ClearFormats()
$ws.Range("H211").NumberFormat = "@"
$ws.Range("H211").Value = "31900.00"
$ws.Range("H211").ClearFormats()
$ws.Range("H212").NumberFormat = "@"
$ws.Range("H212").Value = "2540.94"
$ws.Range("H212").ClearFormats()
$ws.Range("H213").NumberFormat = "@"
$ws.Range("H213").Value = "1200.00"
$ws.Range("H213").ClearFormats()
$ws.Range("H214").NumberFormat = "@"
$ws.Range("H214").Value = "198.00"
$ws.Range("H214").ClearFormats()
